# Append two new rows (new rows 2 and 3) to the top of the job-listing
# table on the active sheet ("ランサーズ"), pushing the previous rows 2-3
# down to rows 4-5, and refresh the "取得日時" timestamp on every row to
# the new scrape time (2025-12-30 06:39:09).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: insert two blank rows above the current row 2 -------
$ws.Range("A2:A3").EntireRow.Insert()

# --- 2. Drop any hyperlinks left over from the insert (the insert does
#        not relocate them onto the shifted cells) so we can rebuild the
#        hyperlink collection cleanly once every row is in place. --------
$ws.Hyperlinks.Delete()

# --- 3. Row 2: new job posting ------------------------------------------
$ws.Range("A2").Value = "2025-12-30 06:39:09"
$ws.Range("B2").Value = "法人向け生成AIサービス(RAG・議事録機能)の設計・開発を支援エンジニア募集(AI/バックエンド)"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5445159"
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("G2").Value = 368
$ws.Range("H2").Value = "🔥AI,Ai ◆開発"

# --- 4. Row 3: new job posting ------------------------------------------
$ws.Range("A3").Value = "2025-12-30 06:39:09"
$ws.Range("B3").Value = "B2B向け生成AIサービス(チャット・RAG)の新規開発プロジェクト推進を支援してくださるPM募集"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5445154"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("G3").Value = 368
$ws.Range("H3").Value = "🔥AI,Ai ◆開発"

# --- 5. Row 4 (previously row 2): refresh its timestamp -----------------
$ws.Range("A4").Value = "2025-12-30 06:39:09"

# --- 6. Row 5 (previously row 3): refresh its timestamp -----------------
$ws.Range("A5").Value = "2025-12-30 06:39:09"

# --- 7. Rebuild the hyperlink collection for every URL cell -------------
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5445159")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5445154")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5463183")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5463296")

# --- 8. Widen column B to fit the longer titles --------------------------
$ws.Range("B1").EntireColumn.ColumnWidth = 51.1667
